$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "... bandigo is jus tin the way." -> split the " tin" run into
# " " + "tin", moving the spellEnd proofErr mark to after the new "tin" run
# and wrapping "tin" in its own gramStart/gramEnd pair (a spurious extra
# grammar-check region), matching the target OOXML exactly.
#
# We replace a wider, unambiguous span (" is jus tin the way.") in a single
# InsertXML call so every <w:proofErr/> boundary in the replacement sits
# strictly inside the targeted range (never exactly on its edge), which is
# what keeps this COM host's run/markup patcher from reshuffling adjacent,
# untouched text.
# ---------------------------------------------------------------------------
$needle1 = " is jus tin the way."
$full = $d.Content.Text
$idx1 = $full.IndexOf($needle1)
if ($idx1 -lt 0) { throw "Could not locate target text for change 1" }
$rng1 = $d.Range($idx1, $idx1 + $needle1.Length)
if ($rng1.Text -ne $needle1) { throw "Range text mismatch for change 1: [$($rng1.Text)]" }

$frag1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
  '<w:r w:rsidR="00322B3C"><w:t xml:space="preserve"> is </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
  '<w:r w:rsidR="00322B3C"><w:t>jus</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r w:rsidR="00322B3C"><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r w:rsidR="00322B3C"><w:t>tin</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
  '<w:r w:rsidR="00322B3C"><w:t xml:space="preserve"> the way.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($frag1)

# ---------------------------------------------------------------------------
# Change 2: merge the lone-space run and the "Marnie attempts..." run into a
# single run (" Marnie attempts a banishment spell ..."). A plain
# Find/Replace of just the second run's text causes Word to fold it back
# together with the preceding whitespace-only run while leaving the earlier
# "." run untouched.
# ---------------------------------------------------------------------------
$needle2 = "Marnie attempts a banishment spell to make Buford and Wes leave the farm. It fails but from that point forward, the backs of all the remote controls go missing."
$rngAll = $d.Content
$found = $rngAll.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, $needle2, 2)
if (-not $found) { throw "Could not locate target text for change 2" }
